$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cls = $m.CustomLayouts
$blank = $cls.Item(7)
$s2 = $p.Slides.AddSlide(2, $blank)

$conn = $s2.Shapes.AddConnector(1, 100, 100, 100, 200)
try {
  $conn.Line.ForeColor.SchemeColor = 2
  Write-Host "SchemeColor set ok:" $conn.Line.ForeColor.SchemeColor
} catch {
  Write-Host "Error:" $_
}
